# Jean's suggestion on Human in the loop for consideration
# Adds two new selected participants to the bottom of the list and
# tidies up the sheet view/column widths to match the reviewed layout.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Column widths (Full name / Gender / Age Group / Email address) ---
$ws.Columns.Item(1).ColumnWidth = 20.830729166666668
$ws.Columns.Item(2).ColumnWidth = 21.385416666666668
$ws.Columns.Item(3).ColumnWidth = 21.276041666666668
$ws.Columns.Item(5).ColumnWidth = 52.498697916666664

# --- New row 36: Dr. John Williams ---
$ws.Cells.Item(2, 1).Copy()
$ws.Cells.Item(36, 1).PasteSpecial(-4122)
$ws.Cells.Item(36, 1).Value = 45418.502638888902
$ws.Cells.Item(36, 2).Value = "Dr. John Williams"
$ws.Cells.Item(36, 3).Value = "Male"
$ws.Cells.Item(36, 4).Value = "21-30"
$ws.Cells.Item(36, 5).Value = "test@abc.com.sg"
$ws.Cells.Item(36, 6).Value = 30338111

# --- New row 37: Candice Washington ---
$ws.Cells.Item(2, 1).Copy()
$ws.Cells.Item(37, 1).PasteSpecial(-4122)
$ws.Cells.Item(37, 1).Value = 45419.1320949074
$ws.Cells.Item(37, 2).Value = "Candice Washington"
$ws.Cells.Item(37, 3).Value = "Female"
$ws.Cells.Item(37, 4).Value = "41-50"
$ws.Cells.Item(37, 5).Value = "test@abc.com.sg"
$ws.Cells.Item(37, 6).Value = 87857672

$excel.CutCopyMode = $false

# --- Update the saved selection to match the reviewer's last position ---
$ws.Range("L28").Select() | Out-Null
